$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new D-column value would be mis-typed as a number by plain
# assignment (e.g. "1.19") are written via a temporary Text number format
# so they stay strings, then the style is reset back to Normal so no stray
# cell-level formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "58.112.31"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "3.115.95"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "527.35"
$ws.Range("E5").Value = "  +0.83%  "
Set-TextValue $ws.Range("D6") "142.34"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.113.51"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "3.651.99"
$ws.Range("E13").Value = "  +0.68%  "
Set-TextValue $ws.Range("D15") "25.82"
$ws.Range("E15").Value = "  -3.70%  "
Set-TextValue $ws.Range("D16") "0.0000165"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "58.130.37"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "3.122.98"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("E19").Value = "  -0.23%  "
Set-TextValue $ws.Range("D20") "12.76"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  -1.33%  "
Set-TextValue $ws.Range("D22") "343.15"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +1.82%  "
Set-TextValue $ws.Range("D25") "67.68"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("E26").Value = "  -0.92%  "
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  +0.62%  "
Set-TextValue $ws.Range("D30") "6.38"
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("E32").Value = "  +1.92%  "
Set-TextValue $ws.Range("D33") "21.07"
$ws.Range("E33").Value = "  +0.42%  "
Set-TextValue $ws.Range("D34") "1.19"
$ws.Range("E34").Value = "  -0.98%  "
Set-TextValue $ws.Range("D35") "158.42"
$ws.Range("E36").Value = "  +0.38%  "
Set-TextValue $ws.Range("D37") "6.19"
$ws.Range("E37").Value = "  +0.89%  "
Set-TextValue $ws.Range("D38") "26.40"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("E39").Value = "  -3.19%  "
Set-TextValue $ws.Range("D40") "1.64"
$ws.Range("E40").Value = "  +12.29%  "
$ws.Range("E41").Value = "  -2.18%  "
Set-TextValue $ws.Range("D42") "4.01"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("D44").Value = "3.156.81"
$ws.Range("E44").Value = "  +0.68%  "
Set-TextValue $ws.Range("D45") "36.54"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D48").Value = "2.271.32"
$ws.Range("E48").Value = "  -0.66%  "
Set-TextValue $ws.Range("D49") "0.998"
$ws.Range("E49").Value = "  +3.74%  "
Set-TextValue $ws.Range("D50") "6.12"
$ws.Range("E50").Value = "  +2.13%  "
Set-TextValue $ws.Range("D51") "20.65"
$ws.Range("E51").Value = "  -1.32%  "

